$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

$ws.Range("L2").Value = 0.98
$ws.Range("L3").Value = 1.06
$ws.Range("L4").Value = 1.14
$ws.Range("L5").Value = 0.82
$ws.Range("H6").Value = 12.19
$ws.Range("L6").Value = 0.98
$ws.Range("L7").Value = 0.91
$ws.Range("L8").Value = 1.13
$ws.Range("L9").Value = 0.92
$ws.Range("H10").Value = 9.550000000000001
$ws.Range("H11").Value = 8.550000000000001
$ws.Range("L11").Value = 1.15
$ws.Range("H12").Value = 7.55
$ws.Range("L12").Value = 1.16
$ws.Range("H13").Value = 6.14
$ws.Range("L13").Value = 1.17
$ws.Range("L14").Value = 1.02
$ws.Range("L15").Value = 1.05
$ws.Range("H16").Value = 3.55
$ws.Range("L16").Value = 0.89
$ws.Range("L17").Value = 0.83
